$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.058.15'
$ws.Range('E2').Value = '  +1.69%  '
$ws.Range('D3').Value = '3.175.17'
$ws.Range('E3').Value = '  +3.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('D8').Value = '3.174.54'
$ws.Range('E8').Value = '  +3.73%  '
$ws.Range('E9').Value = '  +3.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.21'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('E12').Value = '  +1.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000271'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +17.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.78%  '
$ws.Range('D15').Value = '3.695.41'
$ws.Range('E15').Value = '  +3.68%  '
$ws.Range('D16').Value = '65.139.33'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('D17').Value = '3.174.13'
$ws.Range('E17').Value = '  +3.71%  '
$ws.Range('E18').Value = '  +5.23%  '
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '513.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.47%  '
$ws.Range('E21').Value = '  +5.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.732'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.53%  '
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.03%  '
$ws.Range('E28').Value = '  +4.31%  '
$ws.Range('E29').Value = '  +7.18%  '
$ws.Range('E30').Value = '  +6.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +13.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.32'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.75'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0901'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '476.80'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.14'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0421'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.67'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.20%  '
$ws.Range('D42').Value = '3.064.91'
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.286'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.09%  '
$ws.Range('D47').Value = '0.0₃0619'
$ws.Range('E47').Value = '  +19.56%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('E50').Value = '  +6.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.15%  '
